$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.630176381257741
$ws.Range("C2").Value = 0.3335846736778763
$ws.Range("D2").Value = 0.01560609928320034
$ws.Range("F2").Value = 3.529065049541259
$ws.Range("G2").Value = 0.002598343087747254
$ws.Range("I2").Value = 2.060660554140718
$ws.Range("J2").Value = 0.1159874134268719
$ws.Range("L2").Value = 0.397536397468464
$ws.Range("M2").Value = 0.5697315880403409
$ws.Range("N2").Value = 2.484920273678362
$ws.Range("B3").Value = 2.539138205202164
$ws.Range("C3").Value = 0.3049758128916551
$ws.Range("D3").Value = 0.01485929724968571
$ws.Range("F3").Value = 3.519617957760076
$ws.Range("G3").Value = 0.002603648974876902
$ws.Range("I3").Value = 2.062416726719675
$ws.Range("J3").Value = 0.115771194462079
$ws.Range("L3").Value = 0.3957501985477094
$ws.Range("M3").Value = 0.5564649643287751
$ws.Range("N3").Value = 2.508305539574735
$ws.Range("B4").Value = 2.484680832143567
$ws.Range("C4").Value = 0.2875857144509268
$ws.Range("D4").Value = 0.01439417328979076
$ws.Range("F4").Value = 3.51560350650837
$ws.Range("G4").Value = 0.002607079826359898
$ws.Range("I4").Value = 2.064484588275803
$ws.Range("J4").Value = 0.1156377673379527
$ws.Range("L4").Value = 0.3948137588373868
$ws.Range("M4").Value = 0.548615748157161
$ws.Range("N4").Value = 2.523430902619801
$ws.Range("B5").Value = 2.462850947583945
$ws.Range("C5").Value = 0.280542859447678
$ws.Range("D5").Value = 0.01420292241167687
$ws.Range("F5").Value = 3.514415952232113
$ws.Range("G5").Value = 0.002608521577942948
$ws.Range("I5").Value = 2.065575679718435
$ws.Range("J5").Value = 0.1155832247357793
$ws.Range("L5").Value = 0.3944725066374275
$ws.Range("M5").Value = 0.5454916930575635
$ws.Range("N5").Value = 2.529787189233481
$ws.Range("B6").Value = 2.459247964681936
$ws.Range("C6").Value = 0.279376027622277
$ws.Range("D6").Value = 0.014171060567989
$ws.Range("F6").Value = 3.514245818031242
$ws.Range("G6").Value = 0.002608763620582774
$ws.Range("I6").Value = 2.065771846706973
$ws.Range("J6").Value = 0.1155741576545708
$ws.Range("L6").Value = 0.3944182807143264
$ws.Range("M6").Value = 0.5449774495694371
$ws.Range("N6").Value = 2.530854270732121
$ws.Range("B7").Value = 2.484384961573994
$ws.Range("C7").Value = 0.2874905555109422
$ws.Range("D7").Value = 0.01439160100324699
$ws.Range("F7").Value = 3.515585676341814
$ws.Range("G7").Value = 0.002607099093158038
$ws.Range("I7").Value = 2.064498297829644
$ws.Range("J7").Value = 0.1156370324475873
$ws.Range("L7").Value = 0.3948089931323935
$ws.Range("M7").Value = 0.5485733141546802
$ws.Range("N7").Value = 2.523515846288994
$ws.Range("B8").Value = 2.598487250038602
$ws.Range("C8").Value = 0.3236835460194811
$ws.Range("D8").Value = 0.01534993622320968
$ws.Range("F8").Value = 3.525436512545738
$ws.Range("G8").Value = 0.002600136736434681
$ws.Range("I8").Value = 2.061060443515913
$ws.Range("J8").Value = 0.1159129984274667
$ws.Range("L8").Value = 0.3968872737748228
$ws.Range("M8").Value = 0.5650956989636953
$ws.Range("N8").Value = 2.492824251720265
$ws.Range("B9").Value = 2.833699053569944
$ws.Range("C9").Value = 0.3960784024749842
$ws.Range("D9").Value = 0.01717944709796981
$ws.Range("F9").Value = 3.558966121785701
$ws.Range("G9").Value = 0.002587849706623624
$ws.Range("I9").Value = 2.06219152866737
$ws.Range("J9").Value = 0.1164489870159855
$ws.Range("L9").Value = 0.4022332314723229
$ws.Range("M9").Value = 0.5998518065938825
$ws.Range("N9").Value = 2.438730654386987
$ws.Range("B10").Value = 3.013556516148583
$ws.Range("C10").Value = 0.4501734633243473
$ws.Range("D10").Value = 0.01849673190222134
$ws.Range("F10").Value = 3.592328609776217
$ws.Range("G10").Value = 0.002579645918456386
$ws.Range("I10").Value = 2.067854345300049
$ws.Range("J10").Value = 0.1168398234481014
$ws.Range("L10").Value = 0.4069345382729068
$ws.Range("M10").Value = 0.6268309412372517
$ws.Range("N10").Value = 2.402713458065953
$ws.Range("B11").Value = 3.096923911820511
$ws.Range("C11").Value = 0.4749887377130335
$ws.Range("D11").Value = 0.01909092931523304
$ws.Range("F11").Value = 3.609415983001043
$ws.Range("G11").Value = 0.002576090624629944
$ws.Range("I11").Value = 2.071487034434725
$ws.Range("J11").Value = 0.1170170305839093
$ws.Range("L11").Value = 0.4092411801863562
$ws.Range("M11").Value = 0.6394199250781298
$ws.Range("N11").Value = 2.387139709201357
$ws.Range("B12").Value = 3.128716582387824
$ws.Range("C12").Value = 0.484416058560555
$ws.Range("D12").Value = 0.01931527214817308
$ws.Range("F12").Value = 3.616162343462548
$ws.Range("G12").Value = 0.002574769577126681
$ws.Range("I12").Value = 2.073015143576853
$ws.Range("J12").Value = 0.1170840530942545
$ws.Range("L12").Value = 0.4101387791174176
$ws.Range("M12").Value = 0.6442325837424079
$ws.Range("N12").Value = 2.381359169873853
$ws.Range("B13").Value = 3.121859525402101
$ws.Range("C13").Value = 0.4823843596180382
$ws.Range("D13").Value = 0.01926698463035237
$ws.Range("F13").Value = 3.614697114148925
$ws.Range("G13").Value = 0.002575052966723906
$ws.Range("I13").Value = 2.072679247182862
$ws.Range("J13").Value = 0.1170696222063619
$ws.Range("L13").Value = 0.4099443927369464
$ws.Range("M13").Value = 0.6431940680388948
$ws.Range("N13").Value = 2.382598907882951
$ws.Range("B14").Value = 3.099535037048781
$ws.Range("C14").Value = 0.4757637184750934
$ws.Range("D14").Value = 0.019109399163586
$ws.Range("F14").Value = 3.60996547804092
$ws.Range("G14").Value = 0.002575981435334236
$ws.Range("I14").Value = 2.071609693643069
$ws.Range("J14").Value = 0.117022546196166
$ws.Range("L14").Value = 0.4093145429024219
$ws.Range("M14").Value = 0.6398149536094024
$ws.Range("N14").Value = 2.386661795957856
$ws.Range("B15").Value = 3.085889730689587
$ws.Range("C15").Value = 0.4717123482544139
$ws.Range("D15").Value = 0.01901278841098986
$ws.Range("F15").Value = 3.607103154185893
$ws.Range("G15").Value = 0.00257655343569088
$ws.Range("I15").Value = 2.070974436054527
$ws.Range("J15").Value = 0.1169937001733299
$ws.Range("L15").Value = 0.408931882340724
$ws.Range("M15").Value = 0.6377510729695999
$ws.Range("N15").Value = 2.389165666150049
$ws.Range("B16").Value = 3.008139475550593
$ws.Range("C16").Value = 0.4485559401253454
$ws.Range("D16").Value = 0.01845780320685009
$ws.Range("F16").Value = 3.591250441119712
$ws.Range("G16").Value = 0.002579881808233395
$ws.Range("I16").Value = 2.067638249230569
$ws.Range("J16").Value = 0.116828230944229
$ws.Range("L16").Value = 0.4067871713535851
$ws.Range("M16").Value = 0.6260145841663345
$ws.Range("N16").Value = 2.403747579682836
$ws.Range("B17").Value = 2.960839348305683
$ws.Range("C17").Value = 0.4344036639220121
$ws.Range("D17").Value = 0.01811608762022132
$ws.Range("F17").Value = 3.582015373317347
$ws.Range("G17").Value = 0.002581968803646584
$ws.Range("I17").Value = 2.065862623053903
$ws.Range("J17").Value = 0.1167265723929791
$ws.Range("L17").Value = 0.405514463443609
$ws.Range("M17").Value = 0.6188955990220038
$ws.Range("N17").Value = 2.412900972013752
$ws.Range("B18").Value = 2.93377936782241
$ws.Range("C18").Value = 0.4262831159224447
$ws.Range("D18").Value = 0.01791906375777685
$ws.Range("F18").Value = 3.576883381334284
$ws.Range("G18").Value = 0.002583185823876456
$ws.Range("I18").Value = 2.064940753196026
$ws.Range("J18").Value = 0.1166680457527995
$ws.Range("L18").Value = 0.4047982498031075
$ws.Range("M18").Value = 0.6148306917110986
$ws.Range("N18").Value = 2.418242043641079
$ws.Range("B19").Value = 2.92464235156632
$ws.Range("C19").Value = 0.4235369654831516
$ws.Range("D19").Value = 0.01785227098453035
$ws.Range("F19").Value = 3.575176623348923
$ws.Range("G19").Value = 0.002583600747100682
$ws.Range("I19").Value = 2.064645683520183
$ws.Range("J19").Value = 0.1166482200656933
$ws.Range("L19").Value = 0.4045584693398894
$ws.Range("M19").Value = 0.6134594915486034
$ws.Range("N19").Value = 2.42006353034752
$ws.Range("B20").Value = 2.965859434508616
$ws.Range("C20").Value = 0.4359081793357973
$ws.Range("D20").Value = 0.01815251288735098
$ws.Range("F20").Value = 3.582979848928147
$ws.Range("G20").Value = 0.002581744919113487
$ws.Range("I20").Value = 2.066041347680525
$ws.Range("J20").Value = 0.1167373998291588
$ws.Range("L20").Value = 0.4056483089045884
$ws.Range("M20").Value = 0.6196503484345541
$ws.Range("N20").Value = 2.411918680504314
$ws.Range("B21").Value = 3.106086217797952
$ws.Range("C21").Value = 0.4777075342971671
$ws.Range("D21").Value = 0.01915570346662321
$ws.Range("F21").Value = 3.611347782609556
$ws.Range("G21").Value = 0.002575708036917772
$ws.Range("I21").Value = 2.071919704559278
$ws.Range("J21").Value = 0.117036375767567
$ws.Range("L21").Value = 0.4094988906160637
$ws.Range("M21").Value = 0.6408062466900191
$ws.Range("N21").Value = 2.385465251905437
$ws.Range("B22").Value = 3.199034014098402
$ws.Range("C22").Value = 0.5052027777725243
$ws.Range("D22").Value = 0.0198074870902829
$ws.Range("F22").Value = 3.631495586680131
$ws.Range("G22").Value = 0.002571909788869506
$ws.Range("I22").Value = 2.076650630566178
$ws.Range("J22").Value = 0.1172312970790816
$ws.Range("L22").Value = 0.4121560653329936
$ws.Range("M22").Value = 0.6548979839571771
$ws.Range("N22").Value = 2.36885795622495
$ws.Range("B23").Value = 3.149306781840892
$ws.Range("C23").Value = 0.4905116822928903
$ws.Range("D23").Value = 0.01945995166841996
$ws.Range("F23").Value = 3.6205948826005
$ws.Range("G23").Value = 0.002573923561355838
$ws.Range("I23").Value = 2.074044113008611
$ws.Range("J23").Value = 0.1171273067414429
$ws.Range("L23").Value = 0.4107250272152498
$ws.Range("M23").Value = 0.6473526850293112
$ws.Range("N23").Value = 2.377659109650246
$ws.Range("B24").Value = 2.963589437754592
$ws.Range("C24").Value = 0.4352279387973681
$ws.Range("D24").Value = 0.01813604679525227
$ws.Range("F24").Value = 3.582543257075912
$ws.Range("G24").Value = 0.002581846083764953
$ws.Range("I24").Value = 2.065960238055425
$ws.Range("J24").Value = 0.1167325050004315
$ws.Range("L24").Value = 0.4055877491418158
$ws.Range("M24").Value = 0.619309039347506
$ws.Range("N24").Value = 2.412362529476255
$ws.Range("B25").Value = 2.768835058707452
$ws.Range("C25").Value = 0.3763374326318853
$ws.Range("D25").Value = 0.01668948630736011
$ws.Range("F25").Value = 3.548367399620915
$ws.Range("G25").Value = 0.002591028388196186
$ws.Range("I25").Value = 2.06103939850945
$ws.Range("J25").Value = 0.1163045216583374
$ws.Range("L25").Value = 0.4006510933892997
$ws.Range("M25").Value = 0.5901963810665265
$ws.Range("N25").Value = 2.452711027518674
